$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.647016048431396
$ws.Range("B1").Value = 2.220998525619507
$ws.Range("C1").Value = 3.201782464981079
$ws.Range("D1").Value = 4.517483234405518
$ws.Range("E1").Value = 0.6190973520278931
